$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update price (D) and volume (E) cells to the latest scraped values.
# Numeric-looking price values are written via a Text-formatted cell first
# (then restored to the default "Normal" style) so Excel stores them as
# exact text instead of re-parsing them into floating point numbers.
$ws.Range('D2').Value = '27.599.99'
$ws.Range('E2').Value = '  -1.72%  '
$ws.Range('D3').Value = '1.666.52'
$ws.Range('E3').Value = '  -3.47%  '
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '214.81'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.88%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.510'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -2.49%  '
$ws.Range('E7').Value = '  -0.06%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '23.62'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -2.19%  '
$ws.Range('E9').Value = '  -1.05%  '
$ws.Range('E11').Value = '  -2.11%  '
$ws.Range('D12').Value = '1.901.62'
$ws.Range('E12').Value = '  -3.59%  '
$ws.Range('D13').Value = '1.665.65'
$ws.Range('E13').Value = '  -3.55%  '
$ws.Range('E14').Value = '  -3.30%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.561'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -0.15%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '66.17'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -1.98%  '
$ws.Range('D17').Value = '27.608.27'
$ws.Range('E17').Value = '  -1.59%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '243.44'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.28%  '
$ws.Range('D19').Value = '0.0₃0729'
$ws.Range('E19').Value = '  -3.60%  '
$ws.Range('E20').Value = '  -4.36%  '
$ws.Range('E21').Value = '  +0.04%  '
$ws.Range('E22').Value = '  -3.47%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '9.28'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -4.28%  '
$ws.Range('E24').Value = '  -4.35%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '146.73'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -1.50%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '7.19'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -4.26%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '16.45'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -1.60%  '
$ws.Range('E28').Value = '  +0.01%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.112'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -2.31%  '
$ws.Range('E30').Value = '  +2.71%  '
$ws.Range('E32').Value = '  -2.68%  '
$ws.Range('D33').Value = '1.467.24'
$ws.Range('E33').Value = '  -1.85%  '
$ws.Range('E34').Value = '  -5.08%  '
$ws.Range('E35').Value = '  -6.11%  '
$ws.Range('E36').Value = '  -1.76%  '
$ws.Range('E37').Value = '  -2.78%  '
$ws.Range('E38').Value = '  -1.48%  '
$ws.Range('E39').Value = '  -6.03%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '69.38'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -1.79%  '
$ws.Range('E41').Value = '  -5.11%  '
$ws.Range('E42').Value = '  -0.06%  '
$ws.Range('B43').Value = 'MXToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.22'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -3.73%  '
$ws.Range('B44').Value = 'FraxShare'
$ws.Range('C44').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '5.41'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -7.24%  '
$ws.Range('D45').Value = '1.809.43'
$ws.Range('E45').Value = '  -3.54%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.787'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -1.55%  '
$ws.Range('E47').Value = '  -4.05%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '89.29'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -1.99%  '
$ws.Range('E49').Value = '  -4.27%  '
$ws.Range('E50').Value = '  -2.38%  '
$ws.Range('B51').Value = 'BitcoinSV'
$ws.Range('C51').Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '41.05'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +14.32%  '
